$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 9842.3799999999992
$ws.Range("B14").Value = 9883.89
$ws.Range("C14").Value = 311.98
$ws.Range("D14").Value = 310.67
$ws.Range("E14").Value = $false
$ws.Range("F14").Value = -0.42
$ws.Range("G14").Value = 42620.766134259262
$ws.Range("G14").NumberFormat = "m/d/yy h:mm"
$ws.Range("H14").Value = $false
